# "added more comments about nonwhere clauses"
#
# Slide 6 figure tweaks:
#  1. Shrink the "Straight Connector 16" line (shape 14) so its extent
#     (cx) goes from 2441575 EMU (192.25 pt) down to 1603375 EMU (126.25 pt).
#  2. Re-set the text of the small caption textboxes ("complaints", "dirty",
#     "truth", "fix", "(a)", "(b)") so PowerPoint no longer keeps a stray
#     trailing <a:endParaRPr> run-properties echo in those paragraphs.
#     Deleting the existing text range and retyping the same text achieves
#     this cleanly, without touching the run-level formatting.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# 1) Resize the connector line (shape index 14, "Straight Connector 16").
$connector = $s.Shapes.Item(14)
$connector.Width = 126.25

# 2) Retype the caption textboxes to drop the leftover endParaRPr.
$captionShapeIndexes = @(15, 16, 17, 20, 21, 22, 23, 24, 25, 26)

foreach ($idx in $captionShapeIndexes) {
    $shp = $s.Shapes.Item($idx)
    $tr = $shp.TextFrame.TextRange
    $text = $tr.Text
    $tr.Delete()
    $shp.TextFrame.TextRange.Text = $text
}
